# Switch_Calculator.xlsx - "M+ M- Mx M/ -- added"
# Mark rows 175-201 (M-minus, M-mul, M-div entries) as "finished" by
# writing an "x" into column F, mirroring the existing "finished" markers
# used elsewhere on the sheet (see e.g. F3, F5, F7, ...).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 175; $row -le 201; $row++) {
    $ws.Cells.Item($row, 6).Value = "x"
}

# Restore the sheet's frozen-pane scroll position / active selection the
# author left behind when saving.
$null = $ws.Range("F204").Select()
